# Generate Report for Handoff
# ---------------------------------------------------------------
# The localization-status report is regenerated: the job moves from
# "In Translation" to "Ready for handoff", a fresh handoff timestamp is
# stamped, and the Status/Date columns are re-autosized to fit the new
# (longer) text - exactly what a re-run of the report generator does.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value     = $newStatus
$dede.Range("C2").Value     = $newStatus

# --- Timestamps refreshed by the new handoff run ----------------------
# de-de handoff datetime (mirrored onto the Overview "Latest HO Xliff
# Generate Date" cell, which shares the same value)
$dede.Range("H2").Value     = "2016-09-03 09:02:08"
$overview.Range("G2").Value = "2016-09-03 09:02:08"

# zh-cn handoff datetime
$zhcn.Range("H2").Value     = "2016-09-03 09:01:59"

# --- Re-autosize the Status / Date columns for the longer text --------
# ColumnWidth only accepts whole-pixel increments (same as real Excel's
# COM automation), so we feed in calibrated "characters" values that
# land on the closest achievable width to the target layout.
$overview.Range("E1").ColumnWidth = 16.33
$overview.Range("F1").ColumnWidth = 16.33
$zhcn.Range("C1").ColumnWidth     = 16.33
$dede.Range("C1").ColumnWidth     = 16.33
